$d = $word.ActiveDocument
$tab = [char]9

# --- "Setup Assembled by" paragraph: Name: {{assembled_by}}  Date: {{assembled_date}}  Signature: {{assembler_signature}} ---
$d.Content.Find.Execute(
    "{{assembled_by}} " + $tab + $tab + " Date: ", $false, $false, $false, $false, $false, $true, 1, $false,
    "{{assembled_by}}Date: ", 2) | Out-Null

$d.Content.Find.Execute(
    "{{assembled_date}}   " + "       " + "  Signature: ", $false, $false, $false, $false, $false, $true, 1, $false,
    "{{assembled_date}}Signature: ", 2) | Out-Null

# --- "Tested by" paragraph: Name: {{tested_by}}  Date: {{tested_date}}  Signature: {{tester_signature}} ---
$d.Content.Find.Execute(
    "{{tested_by}}            " + $tab + $tab + " Date: ", $false, $false, $false, $false, $false, $true, 1, $false,
    "{{tested_by}}Date: ", 2) | Out-Null

$d.Content.Find.Execute(
    "{{tested_date}} " + "         " + "       " + "   Signature:", $false, $false, $false, $false, $false, $true, 1, $false,
    "{{tested_date}}Signature:", 2) | Out-Null

# --- "Approved by" paragraph: Name: {{approved_by}}  Date: {{approved_date}}  Signature: {{approver_signature}} ---
$d.Content.Find.Execute(
    "{{approved_by}}   " + $tab + $tab + " Date:", $false, $false, $false, $false, $false, $true, 1, $false,
    "{{approved_by}}Date:", 2) | Out-Null

$d.Content.Find.Execute(
    "{{approved_date}} " + "      " + "    " + "   Signature:", $false, $false, $false, $false, $false, $true, 1, $false,
    "{{approved_date}}Signature:", 2) | Out-Null

$d.Saved = $false
